$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.842.74'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.083.91'
$ws.Range('E3').Value = '  +0.62%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '234.15'
$ws.Range('E5').Value = '  -0.37%  '
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.93'
$ws.Range('E7').Value = '  +2.86%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  -0.43%  '
$ws.Range('E10').Value = '  +2.00%  '
$ws.Range('E11').Value = '  +3.04%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.392.92'
$ws.Range('E12').Value = '  +0.67%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.73'
$ws.Range('E13').Value = '  +2.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.22'
$ws.Range('E14').Value = '  +1.88%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.771'
$ws.Range('E15').Value = '  -0.88%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.28'
$ws.Range('E16').Value = '  +2.10%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.091.61'
$ws.Range('E17').Value = '  +0.97%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.805.14'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.19'
$ws.Range('E19').Value = '  +0.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.32'
$ws.Range('E20').Value = '  +2.47%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0833'
$ws.Range('E21').Value = '  +1.80%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '229.08'
$ws.Range('E22').Value = '  +0.87%  '
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('E24').Value = '  -1.15%  '
$ws.Range('E25').Value = '  -1.72%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '170.32'
$ws.Range('E26').Value = '  +1.57%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.139'
$ws.Range('E27').Value = '  +7.56%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.99'
$ws.Range('E28').Value = '  +1.21%  '
$ws.Range('E29').Value = '  +0.26%  '
$ws.Range('E30').Value = '  +2.22%  '
$ws.Range('E32').Value = '  +3.85%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.72'
$ws.Range('E33').Value = '  +4.05%  '
$ws.Range('E34').Value = '  +2.04%  '
$ws.Range('E35').Value = '  +1.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.48'
$ws.Range('E36').Value = '  +2.66%  '
$ws.Range('E37').Value = '  +2.45%  '
$ws.Range('E38').Value = '  -0.05%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.43'
$ws.Range('E39').Value = '  -3.26%  '
$ws.Range('E40').Value = '  +4.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.95'
$ws.Range('E41').Value = '  -0.26%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.47'
$ws.Range('E42').Value = '  +8.39%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '98.58'
$ws.Range('E43').Value = '  +1.66%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0214'
$ws.Range('E44').Value = '  +0.60%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.461.93'
$ws.Range('E45').Value = '  -1.60%  '
$ws.Range('E46').Value = '  +0.83%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '16.06'
$ws.Range('E47').Value = '  +5.35%  '
$ws.Range('E48').Value = '  +3.94%  '
$ws.Range('E49').Value = '  +2.76%  '
$ws.Range('E50').Value = '  +2.78%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.277.70'
$ws.Range('E51').Value = '  +0.58%  '
